# Scheduled-runner data refresh: update the per-leve market/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the affected rows of
# each crafting-job sheet, matching the latest market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3125.6667
$ws.Range("I62").Value = 3591
$ws.Range("K62").Value = 3591
$ws.Range("M62").Value = -2967

$ws.Range("H65").Value = 3125.6667
$ws.Range("I65").Value = 3591
$ws.Range("K65").Value = 17955
$ws.Range("M65").Value = -14835

$ws.Range("H125").Value = 3482.842
$ws.Range("I125").Value = 3060.1667
$ws.Range("J125").Value = 3677.923
$ws.Range("K125").Value = 27541.5003
$ws.Range("L125").Value = 33101.307
$ws.Range("M125").Value = -25081.5003
$ws.Range("N125").Value = -38021.307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1041.6666
$ws.Range("I45").Value = 1041.6666
$ws.Range("K45").Value = 1041.6666
$ws.Range("M45").Value = -664.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4985.0244
$ws.Range("I56").Value = 4985.0244
$ws.Range("K56").Value = 4985.0244
$ws.Range("M56").Value = -4455.0244

$ws.Range("H75").Value = 14975
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 14975
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 44925
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -46921

$ws.Range("H78").Value = 14975
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 14975
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 134775
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -144759

$ws.Range("H133").Value = 4571.3184
$ws.Range("J133").Value = 6281.385
$ws.Range("L133").Value = 18844.155
$ws.Range("N133").Value = -28964.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 32000
$ws.Range("J62").Value = 32000
$ws.Range("L62").Value = 32000
$ws.Range("N62").Value = -33372

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = $null

$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H65").Value = 32000
$ws.Range("J65").Value = 32000
$ws.Range("L65").Value = 96000
$ws.Range("N65").Value = -102864

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = $null

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H68").Value = 31333.334
$ws.Range("J68").Value = 31333.334
$ws.Range("L68").Value = 31333.334
$ws.Range("N68").Value = -32955.334

$ws.Range("H69").Value = 151563.67
$ws.Range("J69").Value = 151563.67
$ws.Range("L69").Value = 151563.67
$ws.Range("N69").Value = -153061.67

$ws.Range("H71").Value = 31333.334
$ws.Range("J71").Value = 31333.334
$ws.Range("L71").Value = 94000.00199999999
$ws.Range("N71").Value = -102112.002

$ws.Range("H72").Value = 151563.67
$ws.Range("J72").Value = 151563.67
$ws.Range("L72").Value = 454691.01
$ws.Range("N72").Value = -462179.01

$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31872

$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99360

$ws.Range("H82").Value = 29250
$ws.Range("J82").Value = 29250
$ws.Range("L82").Value = 29250
$ws.Range("N82").Value = -30016

$ws.Range("H85").Value = 29250
$ws.Range("J85").Value = 29250
$ws.Range("L85").Value = 29250
$ws.Range("N85").Value = -31902

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

$ws.Range("H88").Value = 33695
$ws.Range("J88").Value = 33695
$ws.Range("L88").Value = 33695
$ws.Range("N88").Value = -34597

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null

$ws.Range("H91").Value = 33695
$ws.Range("J91").Value = 33695
$ws.Range("L91").Value = 33695
$ws.Range("N91").Value = -36815

$ws.Range("H122").Value = 3011.8667
$ws.Range("I122").Value = 2086
$ws.Range("J122").Value = 3822
$ws.Range("K122").Value = 6258
$ws.Range("L122").Value = 11466
$ws.Range("M122").Value = -3808
$ws.Range("N122").Value = -16366

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 12000
$ws.Range("I17").Value = 12000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -11830
$ws.Range("N17").Value = $null

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = $null

$ws.Range("H68").Value = 2123.6365
$ws.Range("I68").Value = 1857.5
$ws.Range("K68").Value = 1857.5
$ws.Range("M68").Value = -1108.5

$ws.Range("H71").Value = 2123.6365
$ws.Range("I71").Value = 1857.5
$ws.Range("K71").Value = 9287.5
$ws.Range("M71").Value = -5543.5

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null

$ws.Range("H117").Value = 59933.332
$ws.Range("J117").Value = 59933.332
$ws.Range("L117").Value = 59933.332
$ws.Range("N117").Value = -69111.33199999999

$ws.Range("H123").Value = 32004.592
$ws.Range("J123").Value = 32004.592
$ws.Range("L123").Value = 32004.592
$ws.Range("N123").Value = -41804.592
